# Generate Report for Archive
# The localization status moved on from "Ready for handoff" to "In Translation"
# everywhere that status is reported (the Overview roll-up plus each
# per-locale detail sheet). Once the text shrinks, the status column is
# re-sized to fit the new (shorter) value.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the string literal on the left of -eq. PowerShell's -eq
        # coerces its right operand to the left operand's type, so
        # "$cell.Value2 -eq 'Ready for handoff'" would silently convert the
        # string to $true for boolean cells (e.g. the "Has metadata" column)
        # and produce false-positive matches.
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# Re-size the status columns so their width reflects the shorter text -
# "In Translation" is narrower than "Ready for handoff" was.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = 12.5

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = 12.5
